# Apply the SampleInput.xlsx edit:
#  - rename "Teacher" sheet to "Prof" and restructure its columns
#    (insert "Non_Course" column, rework "Classes Not Available" -> "Required_Course" /
#     "Non_Course", add "Max Credits" column)
#  - Classes sheet: tidy a couple of comma-separated lists and add Room_Exceptions values
#  - Rooms sheet: convert header labels to underscore form and populate the
#    MWF/TTh Room-Time Exceptions columns

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Teacher" -> "Prof"
# ---------------------------------------------------------------------------
$prof = $wb.Worksheets.Item("Teacher")
$prof.Name = "Prof"

# Clear out the old used range first since the new layout has a different
# column count / blank pattern than the old one.
$prof.Cells.Clear()

$prof.Range("A1").Value = "Teacher"
$prof.Range("B1").Value = "Required_Course"
$prof.Range("C1").Value = "Non_Course"
$prof.Range("D1").Value = "MWF_Unavailable_Times"
$prof.Range("E1").Value = "TTh_Unavailable_Times"
$prof.Range("F1").Value = "Room Exceptions"
$prof.Range("G1").Value = "Min Credits"
$prof.Range("H1").Value = "Max Credits"

$prof.Range("A2").Value = "Marchard"
$prof.Range("B2").Value = "MATH 130"
$prof.Range("C2").Value = "MATH 313,STAT 102"
$prof.Range("D2").Value = "8,9,10"
$prof.Range("E2").Value = "9,10"
$prof.Range("G2").Value = 6
$prof.Range("H2").Value = 10

$prof.Range("A3").Value = "Hurl"
$prof.Range("D3").Value = "12,3,4"
$prof.Range("F3").Value = "VSC 201,VSC 203"
$prof.Range("G3").Value = 6
$prof.Range("H3").Value = 12

# ---------------------------------------------------------------------------
# Sheet "Classes"
# ---------------------------------------------------------------------------
$classes = $wb.Worksheets.Item("Classes")

$classes.Range("C2").Value = "MATH 401,MATH 125"
$classes.Range("F2").Value = "VSC 201,VSC 202"

$classes.Range("D3").ClearContents()

$classes.Range("F4").Value = "VSC 203"

$classes.Range("F6").Value = "BML 201"

# ---------------------------------------------------------------------------
# Sheet "Rooms"
# ---------------------------------------------------------------------------
$rooms = $wb.Worksheets.Item("Rooms")

$rooms.Range("B1").Value = "MWF_Room-Time_Exceptions"
$rooms.Range("C1").Value = "TTh_Room-Time_Exceptions"

$rooms.Range("B2").Value = "9,10"
$rooms.Range("C2").Value = 8

$rooms.Range("B3").Value = 11
$rooms.Range("C3").Value = 11

$rooms.Range("B4").Value = 12
$rooms.Range("C4").Value = "12,1"

$rooms.Range("B5").Value = "1,2,3"

# ---------------------------------------------------------------------------
# View state: active sheet becomes "Prof" (first sheet); its sheetView should
# be the tab-selected one, with selection on C14. "Classes" selection moves
# to G3, and "Rooms" selection moves to C5.
# ---------------------------------------------------------------------------
$prof.Activate()
$prof.Range("C14").Select()

$classes.Range("G3").Select()

$rooms.Range("C5").Select()

$prof.Activate()
